$wb = $excel.ActiveWorkbook

# --- Work on the "table_names" sheet (7th sheet) ---
$ws = $wb.Worksheets.Item(7)
$lo = $ws.ListObjects.Item(1)

# Extend the table ("Table26") by one row at the end; this keeps the table
# definition (ref/autoFilter) properly synced with the new row count.
$lo.ListRows.Add() | Out-Null

# Shift existing rows 8-15 down to rows 9-16 (iterate bottom-up to avoid
# clobbering data before it is copied).
for ($r = 15; $r -ge 8; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# event_repeat now displays as "eN" instead of "N" (to disambiguate from
# the new form_repeat row below).
$ws.Cells.Item(6, 1).Value = "eN"

# Insert the new row: form_repeat will now show as "N" in the custom DT.
$ws.Cells.Item(8, 1).Value = "N"
$ws.Cells.Item(8, 2).Value = "form_repeat"

# --- Restore view state for the other sheets touched during editing ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("B9").Select() | Out-Null

# --- Leave the "table_names" sheet active/selected, as it was the sheet
# being edited ---
$ws.Activate()
$ws.Range("A7").Select() | Out-Null
